# create project structure selenium-cucumber-java
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old mailto hyperlinks (and their "Hyperlink" style/underline look)
$ws.Hyperlinks.Delete()
$ws.Range("A2:A4").Style = "Normal"
$wb.Styles("Hyperlink").Delete()

# Rebuild the table: testCase / Email / Password / ExpectedResult
$data = @(
    @("testCase",                   "Email",                    "Password",      "ExpectedResult"),
    @("Login With Wrong Email",     "fiky@gmail.com",           "Usenobi123#",    "ERROR"),
    @("Login With Wrong Password",  "fiky.anggra@usenobi.com",  "Usenobi123",     "ERROR"),
    @("Login Success",              "fiky.anggra@usenobi.com",  "Usenobi123#",    "Coming Soon")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Drop the leftover 5th row / stray "status/number/a/b/c/ss" content
$ws.Rows.Item(5).ClearContents()

# Column widths to fit the new headers
$ws.Columns.Item(1).ColumnWidth = 24.109375
$ws.Columns.Item(2).ColumnWidth = 22.44140625
$ws.Columns.Item(3).ColumnWidth = 11.6640625
$ws.Columns.Item(4).ColumnWidth = 13.5546875

# Selection moves to D9 in the saved view
$ws.Range("D9").Select() | Out-Null
